# Simulated Wild Card round and logged it
# Adds a new player column ("T.Benjamin") to the Players Yards Data workbook.
# The new column is inserted right before the existing "G.Kittle" column
# (i.e. becomes column Q), pushing G.Kittle / R.Dwelley / C.Woerner one
# column to the right on both the "Rushing" and "Receiving" sheets, and
# seeds the new player's yards-count cell with the same placeholder "n"
# value used for every other player.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Insert a new blank column at Q (before the current G.Kittle column),
    # shifting G.Kittle/R.Dwelley/C.Woerner (and anything else) right by one.
    $ws.Columns("Q:Q").Insert()

    # Header row: new player's name.
    $ws.Range("Q1").Value = "T.Benjamin"

    # Data row: same placeholder value ("n") used across the rest of row 2.
    $ws.Range("Q2").Value = "n"
}
